$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New ORM record row (row 8)
$ws.Range("B8").Value = 22041807
$ws.Range("C8").Value = "HLD.docx"
$ws.Range("D8").Value = "Treza Bawn Win"

# Copy the date formatting from an existing "Issue Date" cell so the new
# cell reuses the same style (numFmtId 14 date format) instead of minting
# a new one, then set the date value (43212 = 2018-04-22).
$ws.Range("E2").Copy()
$ws.Range("E8").PasteSpecial(-4122)
$ws.Range("E8").Value = 43212

$excel.CutCopyMode = $false

# Restore the selection to where the user last clicked.
$ws.Range("J12").Select()
